# correção nos dados e inicio da analise PNAD 2009
#
# The sheet was exported from pandas with a couple of spurious
# "header" artifacts:
#   - B2 held the literal text "unnamed: 1_level_1" instead of the
#     real "total" label used by the matching cell in row 1 (B1).
#   - Two rows (originally rows 5 and 8) contained nothing but a
#     leftover pandas multi-index section title - "situação do
#     domicílio" and "grandes regiões e unidades da federação" -
#     with no data of their own; the real data rows sit right below
#     each of them.
#
# Fix: correct the B2 label, then drop those two label-only rows so
# every remaining row carries data, shifting everything else up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled header cell (row 2, column B): "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# Remove the two section-header-only rows. Delete the lower one first
# so the earlier row index (5) still points at the right row when it
# is removed second.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
